# Add a new "JavaScript" worksheet after "Framework", populate it with the
# new paid JavaScript course entry, and make it the active sheet (mirroring
# how "Framework" was previously laid out / selected).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet ("Framework").
$frameworkSheet = $wb.Worksheets.Item("Framework")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $frameworkSheet)
$newSheet.Name = "JavaScript"

# Enter the URL first, then the course title, so the shared-strings table
# ends up ordered [url, title] - matching the source workbook.
$newSheet.Range("C2").Value = "https://www.udemy.com/course/ultimate-javascript-leetcode-interview-bootcamp/"
$newSheet.Range("B2").Value = "JavaScript & LeetCode | The Interview Bootcamp"

# Match the column widths used on the other course sheets.
$newSheet.Columns.Item(2).ColumnWidth = 84.6
$newSheet.Columns.Item(3).ColumnWidth = 99.6

# Leave selection on the next empty row, as on the other sheets.
$newSheet.Range("C4").Select() | Out-Null

# Make the new sheet the active tab.
$newSheet.Activate()
